$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 382.37707044870626
$ws.Range("K6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("T6").Value = 414.9203932492665
$ws.Range("W6").Value = 0
$ws.Range("E7").Value = 1.237075856212648
$ws.Range("K7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("T7").Value = 1.2135839075048422
$ws.Range("W7").Value = 0
$ws.Range("E8").Value = 313.6753901829532
$ws.Range("K8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("T8").Value = 216.5685455702928
$ws.Range("W8").Value = 0
$ws.Range("E9").Value = -0.04461105155872805
$ws.Range("K9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("Q9").Value = 0
$ws.Range("T9").Value = 0.004280489850322706
$ws.Range("W9").Value = 0
$ws.Range("E10").Value = 0.962066901193566
$ws.Range("K10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("T10").Value = 1.008364983618557
$ws.Range("W10").Value = 0
$ws.Range("E11").Value = 0.00010869907345273367
$ws.Range("K11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("T11").Value = -0.000013092854322313775
$ws.Range("W11").Value = 0
$ws.Range("E12").Value = -0.00002234383214367641
$ws.Range("K12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("T12").Value = -0.0001029472370537567
$ws.Range("W12").Value = 0
$ws.Range("E13").Value = -0.00000006574658817802249
$ws.Range("K13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("T13").Value = 0.000000009334881243177565
$ws.Range("W13").Value = 0
$ws.Range("E14").Value = 1.0440689458539125
$ws.Range("K14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("Q14").Value = 0
$ws.Range("T14").Value = 1.0247317426177331
$ws.Range("W14").Value = 0
$ws.Range("E15").Value = -0.00011620582434981155
$ws.Range("K15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("T15").Value = -0.00006502743257830157
$ws.Range("W15").Value = 0
$ws.Range("E16").Value = -0.3269621040302365
$ws.Range("K16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("Q16").Value = 0
$ws.Range("T16").Value = -0.3829840657569394
$ws.Range("W16").Value = 0
$ws.Range("E17").Value = 0.0000000754855217206695
$ws.Range("K17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("T17").Value = 0.00000004194153747272009
$ws.Range("W17").Value = 0
$ws.Range("E18").Value = 0.00026650796299316755
$ws.Range("K18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("T18").Value = 0.00030143555784516277
$ws.Range("W18").Value = 0
$ws.Range("E19").Value = -0.44361520441837204
$ws.Range("K19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("Q19").Value = 0
$ws.Range("T19").Value = -0.41594987678842327
$ws.Range("W19").Value = 0
